$d = $word.ActiveDocument

# Locate the "heading 1" styled paragraph (aliased "כותרת ראשית") -
# the document's single top-level heading paragraph - and explicitly
# set its paragraph-level reading order to left-to-right. This is the
# COM-model equivalent of adding <w:bidi w:val="0"/> to that
# paragraph's <w:pPr> in the OOXML.
foreach ($p in $d.Paragraphs) {
    $styleName = $p.Range.ParagraphFormat.Style.NameLocal
    if ($styleName -like "Heading 1*") {
        $p.Range.ParagraphFormat.ReadingOrder = 0
    }
}
